$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to Text format so numeric-looking strings
# (e.g. "256.77", "98.470.97") are preserved exactly as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '98.470.97'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '3.361.90'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '256.77'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = '665.48'
$ws.Range("E6").Value = '  +6.59%  '
$ws.Range("D7").Value = '1.54'
$ws.Range("E7").Value = '  +7.85%  '
$ws.Range("D8").Value = '0.471'
$ws.Range("E8").Value = '  +20.97%  '
$ws.Range("D9").Value = '1.07'
$ws.Range("E9").Value = '  +21.80%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '3.358.42'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("E12").Value = '  +8.72%  '
$ws.Range("D13").Value = '42.33'
$ws.Range("E13").Value = '  +13.21%  '
$ws.Range("D14").Value = '0.0000274'
$ws.Range("E14").Value = '  +10.62%  '
$ws.Range("D15").Value = '99.613.49'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("D17").Value = '3.982.90'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").Value = '3.358.39'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '7.72'
$ws.Range("E19").Value = '  +26.53%  '
$ws.Range("D20").Value = '16.78'
$ws.Range("E20").Value = '  +10.89%  '
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").Value = '531.56'
$ws.Range("E22").Value = '  +8.65%  '
$ws.Range("E23").Value = '  +13.15%  '
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("D25").Value = '0.435'
$ws.Range("E25").Value = '  +49.16%  '
$ws.Range("D26").Value = '102.42'
$ws.Range("E26").Value = '  +15.57%  '
$ws.Range("D27").Value = '6.22'
$ws.Range("E27").Value = '  +10.88%  '
$ws.Range("D28").Value = '12.57'
$ws.Range("E28").Value = '  +6.42%  '
$ws.Range("D29").Value = '3.534.98'
$ws.Range("E29").Value = '  +1.04%  '
$ws.Range("E30").Value = '  +7.65%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").Value = '11.08'
$ws.Range("E32").Value = '  +14.82%  '
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '29.49'
$ws.Range("E35").Value = '  +5.34%  '
$ws.Range("D36").Value = '0.538'
$ws.Range("E36").Value = '  +17.70%  '
$ws.Range("D37").Value = '7.80'
$ws.Range("E37").Value = '  +7.62%  '
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  +8.76%  '
$ws.Range("D39").Value = '0.157'
$ws.Range("E39").Value = '  +5.20%  '
$ws.Range("D40").Value = '526.47'
$ws.Range("E40").Value = '  +5.86%  '
$ws.Range("E41").Value = '  +6.46%  '
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").Value = '3.93'
$ws.Range("E42").Value = '  +7.01%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").Value = '24.70'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '0.0436'
$ws.Range("E44").Value = '  +34.65%  '
$ws.Range("E45").Value = '  +3.98%  '
$ws.Range("D46").Value = '0.826'
$ws.Range("E46").Value = '  +5.71%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  +6.76%  '
$ws.Range("D49").Value = '1.55'
$ws.Range("E49").Value = '  +13.15%  '
$ws.Range("D50").Value = '5.12'
$ws.Range("E50").Value = '  +11.17%  '
$ws.Range("D51").Value = '50.96'
$ws.Range("E51").Value = '  +11.65%  '
